# Refresh the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures on Sheet1 with the latest scraped snapshot.
#
# Some of the new Price strings look like plain numbers to Excel (e.g.
# "1.004", "0.07670"), which would otherwise be auto-converted from text
# to a numeric value on assignment (losing the original formatting, such
# as trailing zeros). Prefixing those particular values with a leading
# apostrophe forces Excel to keep them as literal text, matching how the
# source data is stored. Values that are unambiguous as text (e.g.
# "26.370.74", the percentage strings in column E) are written directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.370.74'
$ws.Range("E2").Value = '  +1.30%  '
$ws.Range("D3").Value = '1.843.25'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D4").Value = "'" + '1.004'
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").Value = "'" + '259.27'
$ws.Range("E5").Value = '  -6.94%  '
$ws.Range("D6").Value = "'" + '1.003'
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("D7").Value = "'" + '0.5193'
$ws.Range("E7").Value = '  +1.88%  '
$ws.Range("D8").Value = "'" + '0.3265'
$ws.Range("E8").Value = '  -6.65%  '
$ws.Range("D9").Value = "'" + '0.06736'
$ws.Range("E9").Value = '  -1.22%  '
$ws.Range("D10").Value = "'" + '19.35'
$ws.Range("E10").Value = '  -2.74%  '
$ws.Range("D11").Value = "'" + '0.7720'
$ws.Range("E11").Value = '  -4.58%  '
$ws.Range("D12").Value = "'" + '0.07670'
$ws.Range("E12").Value = '  -1.29%  '
$ws.Range("D13").Value = '1.876.20'
$ws.Range("E13").Value = '  +1.97%  '
$ws.Range("D14").Value = "'" + '88.12'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").Value = "'" + '5.047'
$ws.Range("E15").Value = '  -0.65%  '
$ws.Range("D16").Value = "'" + '1.004'
$ws.Range("E16").Value = '  +0.48%  '
$ws.Range("D17").Value = "'" + '14.13'
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("E18").Value = '  +0.28%  '
$ws.Range("D19").Value = "'" + '0.000007896'
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("D20").Value = '26.429.36'
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("D21").Value = '2.061.20'
$ws.Range("E21").Value = '  -0.42%  '
$ws.Range("D22").Value = "'" + '4.601'
$ws.Range("E22").Value = '  -3.68%  '
$ws.Range("D23").Value = "'" + '9.626'
$ws.Range("E23").Value = '  -4.58%  '
$ws.Range("D24").Value = "'" + '5.996'
$ws.Range("E24").Value = '  -3.48%  '
$ws.Range("D25").Value = "'" + '2.323'
$ws.Range("E25").Value = '  -1.67%  '
$ws.Range("D26").Value = "'" + '145.15'
$ws.Range("E26").Value = '  +1.10%  '
$ws.Range("D27").Value = "'" + '1.637'
$ws.Range("E27").Value = '  -1.48%  '
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("D29").Value = "'" + '111.47'
$ws.Range("E29").Value = '  +1.66%  '
$ws.Range("D30").Value = "'" + '4.229'
$ws.Range("E30").Value = '  -2.95%  '
$ws.Range("D31").Value = "'" + '4.184'
$ws.Range("E31").Value = '  -2.41%  '
$ws.Range("D32").Value = "'" + '0.08742'
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("D33").Value = "'" + '0.04850'
$ws.Range("E33").Value = '  -0.21%  '
$ws.Range("D34").Value = "'" + '1.138'
$ws.Range("E34").Value = '  -2.16%  '
$ws.Range("D35").Value = "'" + '2.866'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").Value = "'" + '0.7073'
$ws.Range("E36").Value = '  -3.09%  '
$ws.Range("D37").Value = "'" + '3.102'
$ws.Range("E37").Value = '  -3.49%  '
$ws.Range("D38").Value = "'" + '0.01810'
$ws.Range("E38").Value = '  -2.22%  '
$ws.Range("D39").Value = "'" + '2.216'
$ws.Range("E39").Value = '  -6.58%  '
$ws.Range("D40").Value = "'" + '0.4939'
$ws.Range("E40").Value = '  -4.23%  '
$ws.Range("D41").Value = "'" + '112.52'
$ws.Range("E41").Value = '  -3.82%  '
$ws.Range("D42").Value = "'" + '0.9061'
$ws.Range("E42").Value = '  -4.89%  '
$ws.Range("D43").Value = "'" + '6.087'
$ws.Range("E43").Value = '  -2.76%  '
$ws.Range("D45").Value = "'" + '7.781'
$ws.Range("E45").Value = '  -2.77%  '
$ws.Range("D46").Value = "'" + '0.4275'
$ws.Range("E46").Value = '  -5.44%  '
$ws.Range("D47").Value = "'" + '0.1287'
$ws.Range("E47").Value = '  -5.44%  '
$ws.Range("D48").Value = "'" + '9.234'
$ws.Range("E48").Value = '  -0.54%  '
$ws.Range("E49").Value = '  +0.07%  '
$ws.Range("D50").Value = "'" + '35.30'
$ws.Range("E50").Value = '  -2.29%  '
$ws.Range("D51").Value = "'" + '1.438'
$ws.Range("E51").Value = '  -3.74%  '
